# Fix row 774: A774 and D774 were stored as text; convert them to numbers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(774, 1).Value = 27
$ws.Cells.Item(774, 4).Value = 22

# Append new weather rows 775-861
$data = @(
    @(28, "Cloudy", "08/28/2024", 18, 0, 0),
    @(28, "Cloudy", "08/28/2024", 18, 0, 0),
    @(28, "Cloudy", "08/28/2024", 19, 0, 0),
    @(28, "Cloudy", "08/28/2024", 19, 0, 0),
    @(28, "Cloudy", "08/28/2024", 19, 0, 0),
    @(28, "Cloudy", "08/28/2024", 19, 0, 0),
    @(28, "Cloudy", "08/28/2024", 19, 0, 0),
    @(27, "Cloudy", "08/28/2024", 20, 0, 0),
    @(28, "Cloudy", "08/28/2024", 20, 0, 0),
    @(28, "Cloudy", "08/28/2024", 20, 0, 0),
    @(27, "Cloudy", "08/28/2024", 20, 0, 0),
    @(27, "Cloudy", "08/28/2024", 20, 0, 0),
    @(27, "Cloudy", "08/28/2024", 20, 0, 0),
    @(27, "Cloudy", "08/28/2024", 21, 0, 0),
    @(27, "Cloudy", "08/28/2024", 21, 0, 0),
    @(27, "Cloudy", "08/28/2024", 21, 0, 0),
    @(27, "Cloudy", "08/28/2024", 21, 0, 0),
    @(27, "Cloudy", "08/28/2024", 21, 0, 0),
    @(27, "Cloudy", "08/28/2024", 22, 0, 0),
    @(27, "Cloudy", "08/28/2024", 22, 0, 0),
    @(27, "Cloudy", "08/28/2024", 22, 0, 0),
    @(27, "Cloudy", "08/28/2024", 22, 0, 0),
    @(27, "Cloudy", "08/28/2024", 22, 0, 0),
    @(27, "Cloudy", "08/28/2024", 23, 0, 0),
    @(27, "Cloudy", "08/28/2024", 23, 0, 0),
    @(27, "Cloudy", "08/28/2024", 23, 0, 0),
    @(27, "Cloudy", "08/28/2024", 23, 0, 0),
    @(27, "Cloudy", "08/28/2024", 23, 0, 0),
    @(27, "Cloudy", "08/28/2024", 23, 0, 0),
    @(27, "Cloudy", "08/29/2024", 0, 0, 0),
    @(28, "Cloudy", "08/29/2024", 8, 0, 0),
    @(28, "Cloudy", "08/29/2024", 8, 0, 0),
    @(28, "Cloudy", "08/29/2024", 8, 0, 0),
    @(29, "Light Rain", "08/29/2024", 8, 0, 0),
    @(27, "Rain", "08/29/2024", 8, 0, 0),
    @(29, "Cloudy", "08/29/2024", 8, 0, 0),
    @(29, "Light Rain", "08/29/2024", 9, 0, 0),
    @(29, "Cloudy", "08/29/2024", 9, 0, 0),
    @(29, "Cloudy", "08/29/2024", 9, 0, 0),
    @(29, "Cloudy", "08/29/2024", 9, 0, 0),
    @(29, "Cloudy", "08/29/2024", 9, 0, 0),
    @(29, "Cloudy", "08/29/2024", 9, 0, 0),
    @(29, "Cloudy", "08/29/2024", 10, 0, 0),
    @(29, "Cloudy", "08/29/2024", 10, 0, 0),
    @(29, "Cloudy", "08/29/2024", 10, 0, 0),
    @(30, "Cloudy", "08/29/2024", 10, 0, 0),
    @(30, "Cloudy", "08/29/2024", 10, 0, 0),
    @(30, "Cloudy", "08/29/2024", 11, 0, 0),
    @(30, "Cloudy", "08/29/2024", 11, 0, 0),
    @(30, "Cloudy", "08/29/2024", 11, 0, 0),
    @(30, "Cloudy", "08/29/2024", 11, 0, 0),
    @(30, "Cloudy", "08/29/2024", 11, 0, 0),
    @(29, "Cloudy", "08/29/2024", 11, 0, 0),
    @(30, "Cloudy", "08/29/2024", 12, 0, 0),
    @(29, "Cloudy", "08/29/2024", 12, 0, 0),
    @(29, "Cloudy", "08/29/2024", 12, 0, 0),
    @(30, "Cloudy", "08/29/2024", 12, 0, 0),
    @(30, "Cloudy", "08/29/2024", 12, 0, 0),
    @(30, "Cloudy", "08/29/2024", 12, 0, 0),
    @(30, "Cloudy", "08/29/2024", 12, 0, 0),
    @(30, "Cloudy", "08/29/2024", 13, 0, 0),
    @(30, "Cloudy", "08/29/2024", 13, 0, 0),
    @(30, "Cloudy", "08/29/2024", 13, 0, 0),
    @(30, "Cloudy", "08/29/2024", 13, 0, 0),
    @(30, "Cloudy", "08/29/2024", 13, 0, 0),
    @(29, "Cloudy", "08/29/2024", 13, 0, 0),
    @(29, "Cloudy", "08/29/2024", 14, 0, 0),
    @(29, "Cloudy", "08/29/2024", 14, 0, 0),
    @(29, "Cloudy", "08/29/2024", 14, 0, 0),
    @(29, "Cloudy", "08/29/2024", 14, 0, 0),
    @(29, "Cloudy", "08/29/2024", 14, 0, 0),
    @(29, "Cloudy", "08/29/2024", 14, 0, 0),
    @(29, "Cloudy", "08/29/2024", 15, 0, 0),
    @(29, "Cloudy", "08/29/2024", 15, 0, 0),
    @(29, "Cloudy", "08/29/2024", 15, 0, 0),
    @(29, "Cloudy", "08/29/2024", 15, 0, 0),
    @(29, "Cloudy", "08/29/2024", 16, 0, 0),
    @(29, "Cloudy", "08/29/2024", 16, 0, 0),
    @(29, "Cloudy", "08/29/2024", 16, 0, 0),
    @(28, "Cloudy", "08/29/2024", 16, 0, 0),
    @(29, "Cloudy", "08/29/2024", 16, 0, 0),
    @(28, "Cloudy", "08/29/2024", 16, 0, 0),
    @(28, "Cloudy", "08/29/2024", 17, 0, 0),
    @(28, "Cloudy", "08/29/2024", 17, 0, 0),
    @(28, "Cloudy", "08/29/2024", 17, 0, 0),
    @(28, "Cloudy", "08/29/2024", 17, 0, 0),
    @(28, "Cloudy", "08/29/2024", 17, 1, 1)
)

$startRow = 775
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $data[$i]
    $rowIndex = $startRow + $i
    $aVal = $r[0]
    $bVal = $r[1]
    $cVal = $r[2]
    $dVal = $r[3]
    $aIsText = $r[4]
    $dIsText = $r[5]

    $aCell = $ws.Cells.Item($rowIndex, 1)
    if ($aIsText -eq 1) {
        $aCell.Value = "'" + $aVal
        $aCell.Style = "Normal"
    } else {
        $aCell.Value = [double]$aVal
    }

    $ws.Cells.Item($rowIndex, 2).Value = $bVal

    $cCell = $ws.Cells.Item($rowIndex, 3)
    $cCell.Value = "'" + $cVal
    $cCell.Style = "Normal"

    $dCell = $ws.Cells.Item($rowIndex, 4)
    if ($dIsText -eq 1) {
        $dCell.Value = "'" + $dVal
        $dCell.Style = "Normal"
    } else {
        $dCell.Value = [double]$dVal
    }
}

